$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 161
$ws.Range("I2").Value = 499
$ws.Range("J2").Value = 1996
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 508
$ws.Range("M2").Value = 42
$ws.Range("N2").Value = 341
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 28
$ws.Range("S2").Value = 195
$ws.Range("T2").Value = 346
$ws.Range("U2").Value = 26
$ws.Range("V2").Value = 2963
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 2959
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 57
$ws.Range("AA2").Value = 14
